# Add I0 and IF columns to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns I and J
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Match the header style used by the existing header row (e.g. H1) by
# copying its formatting (not value) onto the new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill data rows 2 through 36: I = 1 (constant), J = same value as column H
for ($r = 2; $r -le 36; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
